$d = $word.ActiveDocument

# Insert a new centered paragraph with contact info right after the
# "Dheeraj Chand" name line. Using Find/Execute with a "^p" in the
# replacement text splits the paragraph at the matched text's end,
# producing a new paragraph that inherits the original paragraph's
# formatting (w:jc center) without inheriting the matched run's
# character formatting (bold / sz 28) onto the new run.
$d.Content.Find.Execute(
    "Dheeraj Chand",
    $true,
    $false,
    $false,
    $false,
    $false,
    $true,
    1,
    $false,
    "Dheeraj Chand^p202.550.7110 | dheeraj.chand@gmail.com | https://www.dheerajchand.com | https://www.linkedin.com/in/dheerajchand/ | Austin, TX",
    2
)
